# Replace the five small inline "legend" pictures in the Waterbodies document
# with plain hyperlink runs pointing at the corresponding image URLs on
# ura.gov.sg, per the commit's intent ("yay the docx works" — the inline
# base64 placeholder pictures get swapped for real hyperlinks to the
# hosted images). The order of the pictures in the document (by their
# AlternativeText) maps to the URLs below.

$d = $word.ActiveDocument

$urls = @(
    "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Commercial/C16_Waterbodies_1.jpg?h=100%25&w=100%25",
    "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Commercial/C17_Waterbodies_2.jpg?h=100%25&w=100%25",
    "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Commercial/C13_Foreshore_A.jpg?h=100%25&w=100%25",
    "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Commercial/C14_Foreshore_B.jpg?h=100%25&w=100%25",
    "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Commercial/C15_Foreshore_C.jpg?h=100%25&w=100%25"
)

$count = $d.InlineShapes.Count
for ($i = 0; $i -lt $count; $i++) {
    # Each picture removed collapses the collection, so the next one to
    # process is always back at index 1.
    $shape = $d.InlineShapes.Item(1)
    $range = $shape.Range
    $url = $urls[$i]

    # Drop the picture itself, leaving an empty range in its place …
    $shape.Delete()

    # … then turn that spot into a hyperlink run whose visible text is the
    # image URL itself (matching the Hyperlink character style).
    $d.Hyperlinks.Add($range, $url, $null, $null, $url) | Out-Null
}
